# Cart sheet: the prices were entered as text (sharedString "$XX.XX").
# Replace each with the plain numeric amount so the column holds real
# numbers instead of label-like strings, and leave the selection where
# the author ended up after editing the last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cart")
$ws.Activate()

$ws.Range("E2").Value = 29.99
$ws.Range("E3").Value = 9.99
$ws.Range("E4").Value = 15.99
$ws.Range("E5").Value = 49.99
$ws.Range("E6").Value = 7.99
$ws.Range("E7").Value = 15.99

$ws.Range("E8").Select()
